# Weekly update: insert a new daily price record at row 7 (pushing the
# existing records from row 7 downward by one row) for the
# "Hortaliza, Vega Central Mapocho de Santiago - Arveja Verde" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 7; Excel shifts rows 7..84 down
# to 8..85 and extends the used range to A1:R85 automatically.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the latest market entry.
$ws.Cells.Item(7, 1).Value = 9
$ws.Cells.Item(7, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(7, 3).Value = "Metropolitana"
$ws.Cells.Item(7, 4).Value = 44530
$ws.Cells.Item(7, 5).Value = 13
$ws.Cells.Item(7, 6).Value = 100112022
$ws.Cells.Item(7, 7).Value = "Arveja Verde"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 43
$ws.Cells.Item(7, 11).Value = 16000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 17023
$ws.Cells.Item(7, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(7, 15).Value = "Carahue"
$ws.Cells.Item(7, 16).Value = 681
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
